$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 0.8160085082054138
$ws.Range("B1").Value = 1.810768961906433
$ws.Range("C1").Value = 6.709448337554932
$ws.Range("D1").Value = 1.597153425216675
$ws.Range("E1").Value = 0.91977858543396
